$d = $word.ActiveDocument

# Clear the "Student Name" value cell (Table 1, Row 1, Col 2)
$table = $d.Tables.Item(1)
$nameCell = $table.Cell(1, 2)
$nameCell.Range.Text = ""

# Clear the "Student No" value cell (Table 1, Row 2, Col 2), including the _GoBack bookmark
$noCell = $table.Cell(2, 2)
$noCell.Range.Text = ""
